# Rasha Comment _ Code Enhancements
# Rename test data method-name strings to use underscores.

$wb = $excel.ActiveWorkbook

$wsCustomer = $wb.Worksheets.Item("Customer Data")
$wsManager  = $wb.Worksheets.Item("Manager Data")

$wsCustomer.Range("A2").Value = "Create_New_User1"
$wsCustomer.Range("A3").Value = "Create_New_User2"

$wsManager.Range("A2").Value = "Manager_Login"

# Leave the last-used selection on "Customer Data" at A3 (matches the
# recorded sheetView state after the edit), then restore "Manager Data"
# as the active sheet/tab so the workbook-level view state is unchanged.
$wsCustomer.Range("A3").Select() | Out-Null
$wsManager.Activate() | Out-Null
